$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Name) and Column B (Link) values for rows 1-15, in final order.
$colA = @(
    "Name",
    "Clustering Billions of Reads for DNA Data Storage",
    "Robust Hierarchical Clustering",
    "EmbedJoin: Efficient Edit Similarity Joins via Embeddings",
    "Low Distortion Embeddings for Edit Distance",
    "A Comprehensive Survey of Clustering Algorithms",
    "Survey of Clustering Algorithms",
    "A Characterization of the DNA Data Storage Channel",
    "Flexible Models for Microclustering with Application to Entity Resolution",
    "A Hierarchical Algorithm for Extreme Clustering",
    "Convolutional Embedding for Edit Distance",
    "Deep metric learning using Triplet network",
    "Unsupervised Deep Embedding for Clustering Analysis",
    "Clustering with Deep Learning: Taxonomy and New Methods",
    "A Survey of Clustering With Deep Learning: From the Perspective of Network Architecture"
)

$colB = @(
    "Link",
    "https://www.microsoft.com/en-us/research/publication/clustering-billions-of-reads-for-dna-data-storage/",
    "https://www.jmlr.org/papers/volume15/balcan14a/balcan14a.pdf",
    "https://dl.acm.org/doi/abs/10.1145/3097983.3098003",
    "https://citeseerx.ist.psu.edu/viewdoc/download?doi=10.1.1.1077.3119&rep=rep1&type=pdf",
    "https://link.springer.com/article/10.1007/s40745-015-0040-1",
    "https://scholarsmine.mst.edu/cgi/viewcontent.cgi?article=1763&context=ele_comeng_facwork",
    "https://www.nature.com/articles/s41598-019-45832-6#Sec7",
    "https://papers.nips.cc/paper/2016/file/670e8a43b246801ca1eaca97b3e19189-Paper.pdf",
    "https://dl.acm.org/doi/abs/10.1145/3097983.3098079",
    "https://arxiv.org/abs/2001.11692",
    "https://arxiv.org/abs/1412.6622",
    "http://proceedings.mlr.press/v48/xieb16.html",
    "https://arxiv.org/abs/1801.07648",
    "https://ieeexplore.ieee.org/abstract/document/8412085"
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
}

# Highlight rows 2, 6, 8 and 15 ("Good" cell style - green fill).
$highlightRows = @(2, 6, 8, 15)
foreach ($r in $highlightRows) {
    $ws.Range("A$r`:B$r").Style = "Good"
}

# Restore the selection state to the whole of row 15, matching the saved view.
$ws.Range("A15:XFD15").Select()
